$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.09665433333333333
$ws.Range("H2").Value = 0.289963
$ws.Range("I2").Value = 0.0006230336790718351
$ws.Range("J2").Value = 0.0006230336790718351
$ws.Range("M2").Value = 4.47329
$ws.Range("N2").Value = 13.41987
$ws.Range("O2").Value = 0.3468876470949054
$ws.Range("P2").Value = 0.3468876470949054
$ws.Range("Q2").Value = 0.4323628627566666
$ws.Range("R2").Value = 3.89126576481
$ws.Range("S2").Value = 0.0002161226869941112
$ws.Range("T2").Value = 0.0002161226869941112
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.09665433333333333
$ws.Range("H3").Value = 0.289963
$ws.Range("I3").Value = 0.0006230336790718351
$ws.Range("J3").Value = 0.0006230336790718351
$ws.Range("O3").Value = 0.3372845821706862
$ws.Range("P3").Value = 0.3372845821706862
$ws.Range("Q3").Value = 0.4203935445158888
$ws.Range("R3").Value = 3.783541900643
$ws.Range("S3").Value = 0.0002101396541240093
$ws.Range("T3").Value = 0.0002101396541240093
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.09665433333333333
$ws.Range("H4").Value = 0.289963
$ws.Range("I4").Value = 0.0006230336790718351
$ws.Range("J4").Value = 0.0006230336790718351
$ws.Range("M4").Value = 4.072757333333333
$ws.Range("N4").Value = 12.218272
$ws.Range("O4").Value = 0.3158277707344083
$ws.Range("P4").Value = 0.3158277707344083
$ws.Range("Q4").Value = 0.3936496448817777
$ws.Range("R4").Value = 3.542846803935999
$ws.Range("S4").Value = 0.0001967713379537144
$ws.Range("T4").Value = 0.0001967713379537144
$ws.Range("G5").Value = 154.8642143333334
$ws.Range("H5").Value = 464.5926430000001
$ws.Range("I5").Value = 0.9982544794956518
$ws.Range("J5").Value = 0.9982544794956519
$ws.Range("M5").Value = 4.47329
$ws.Range("N5").Value = 13.41987
$ws.Range("O5").Value = 0.3468876470949054
$ws.Range("P5").Value = 0.3468876470949054
$ws.Range("Q5").Value = 692.7525413351567
$ws.Range("R5").Value = 6234.77287201641
$ws.Range("S5").Value = 0.3462821475941961
$ws.Range("T5").Value = 0.3462821475941962
$ws.Range("G6").Value = 154.8642143333334
$ws.Range("H6").Value = 464.5926430000001
$ws.Range("I6").Value = 0.9982544794956518
$ws.Range("J6").Value = 0.9982544794956519
$ws.Range("O6").Value = 0.3372845821706862
$ws.Range("P6").Value = 0.3372845821706862
$ws.Range("Q6").Value = 673.5747248675693
$ws.Range("S6").Value = 0.3366958450167068
$ws.Range("T6").Value = 0.3366958450167069
$ws.Range("G7").Value = 154.8642143333334
$ws.Range("H7").Value = 464.5926430000001
$ws.Range("I7").Value = 0.9982544794956518
$ws.Range("J7").Value = 0.9982544794956519
$ws.Range("M7").Value = 4.072757333333333
$ws.Range("N7").Value = 12.218272
$ws.Range("O7").Value = 0.3158277707344083
$ws.Range("P7").Value = 0.3158277707344083
$ws.Range("Q7").Value = 630.7243645969885
$ws.Range("R7").Value = 5676.519281372896
$ws.Range("S7").Value = 0.3152764868847488
$ws.Range("T7").Value = 0.3152764868847488
$ws.Range("G8").Value = 0.174137
$ws.Range("H8").Value = 0.522411
$ws.Range("I8").Value = 0.001122486825276316
$ws.Range("J8").Value = 0.001122486825276316
$ws.Range("M8").Value = 4.47329
$ws.Range("N8").Value = 13.41987
$ws.Range("O8").Value = 0.3468876470949054
$ws.Range("P8").Value = 0.3468876470949054
$ws.Range("Q8").Value = 0.7789653007299998
$ws.Range("R8").Value = 7.010687706569999
$ws.Range("S8").Value = 0.0003893768137151314
$ws.Range("T8").Value = 0.0003893768137151314
$ws.Range("G9").Value = 0.174137
$ws.Range("H9").Value = 0.522411
$ws.Range("I9").Value = 0.001122486825276316
$ws.Range("J9").Value = 0.001122486825276316
$ws.Range("O9").Value = 0.3372845821706862
$ws.Range("P9").Value = 0.3372845821706862
$ws.Range("Q9").Value = 0.7574008131523332
$ws.Range("R9").Value = 6.816607318370999
$ws.Range("S9").Value = 0.0003785974998554223
$ws.Range("T9").Value = 0.0003785974998554223
$ws.Range("G10").Value = 0.174137
$ws.Range("H10").Value = 0.522411
$ws.Range("I10").Value = 0.001122486825276316
$ws.Range("J10").Value = 0.001122486825276316
$ws.Range("M10").Value = 4.072757333333333
$ws.Range("N10").Value = 12.218272
$ws.Range("O10").Value = 0.3158277707344083
$ws.Range("P10").Value = 0.3158277707344083
$ws.Range("Q10").Value = 0.7092177437546665
$ws.Range("R10").Value = 6.382959693791999
$ws.Range("S10").Value = 0.000354512511705762
$ws.Range("T10").Value = 0.0003545125117057622
